{"js": "// 1) Remove the run of 4 spaces that sits right before \"demonstration\"\n//    (the second occurrence, immediately preceded by the red error text).\nconst scoped = context.document.body.search(\n  \"345]    demonstration\",\n  { matchCase: true }\n);\nscoped.load(\"items\");\nawait context.sync();\n\nconst spaces = scoped.items[0].search(\"    \", { matchCase: true });\nspaces.load(\"items\");\nawait context.sync();\n\nspaces.items[0].insertText(\"\", \"Replace\");\nawait context.sync();\n\n// 2) Replace the field-code construct ({ m:self.name }, built from\n//    fldChar/instrText runs) with plain literal text runs:\n//    \"{\" \"m\" \":\" \"self\" \".name}\"  (the \"self\" run keeps its color).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst fieldParagraph = paragraphs.items[1];\nconst fieldRange = fieldParagraph.getRange();\n\nconst replacementOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r><w:t>{</w:t></w:r>\n            <w:r><w:t>m</w:t></w:r>\n            <w:r><w:t>:</w:t></w:r>\n            <w:r>\n              <w:rPr>\n                <w:color w:themeColor=\"accent6\" w:themeShade=\"BF\" w:val=\"E36C0A\"/>\n              </w:rPr>\n              <w:t>self</w:t>\n            </w:r>\n            <w:r><w:t xml:space=\"preserve\">.name}</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\nfieldRange.insertOoxml(replacementOoxml, \"Replace\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Remove the run of 4 spaces that sits right before \"demonstration\"\n#    (the second occurrence, immediately preceded by the red error text).\n$scope = $d.Content\n$scope.Find.ClearFormatting()\n$scope.Find.Execute(\"345]    demonstration\") | Out-Null\n\n$spaces = $d.Range($scope.Start, $scope.End)\n$spaces.Find.ClearFormatting()\n$spaces.Find.Execute(\"    \") | Out-Null\n$spaces.Text = \"\"\n\n# 2) Replace the field-code construct ({ m:self.name }, built from\n#    fldChar/instrText runs) with plain literal text runs:\n#    \"{\" \"m\" \":\" \"self\" \".name}\"  (the \"self\" run keeps its color).\n$fieldParagraph = $d.Paragraphs.Item(2)\n$fieldRange = $fieldParagraph.Range\n\n$replacementOoxml = @\"\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r><w:t>{</w:t></w:r>\n            <w:r><w:t>m</w:t></w:r>\n            <w:r><w:t>:</w:t></w:r>\n            <w:r>\n              <w:rPr>\n                <w:color w:themeColor=\"accent6\" w:themeShade=\"BF\" w:val=\"E36C0A\"/>\n              </w:rPr>\n              <w:t>self</w:t>\n            </w:r>\n            <w:r><w:t xml:space=\"preserve\">.name}</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n\"@\n\n$fieldRange.InsertXML($replacementOoxml)\n"}
